# Apply cell value updates from cryptos list refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '27.678.78'
$c.Style = $origStyle
$ws.Range('E2').Value = '  -0.93%  '
$c = $ws.Range('D3')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.583.14'
$c.Style = $origStyle
$ws.Range('E3').Value = '  -3.22%  '
$ws.Range('E4').Value = '  +0.40%  '
$c = $ws.Range('D5')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '206.49'
$c.Style = $origStyle
$ws.Range('E5').Value = '  -2.44%  '
$ws.Range('E6').Value = '  -2.61%  '
$ws.Range('E7').Value = '  +0.45%  '
$c = $ws.Range('D8')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '22.24'
$c.Style = $origStyle
$ws.Range('E8').Value = '  -4.87%  '
$ws.Range('E9').Value = '  -1.66%  '
$ws.Range('E10').Value = '  -3.20%  '
$c = $ws.Range('D11')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.0867'
$c.Style = $origStyle
$ws.Range('E11').Value = '  -1.81%  '
$c = $ws.Range('D12')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.808.44'
$c.Style = $origStyle
$ws.Range('E12').Value = '  -3.15%  '
$c = $ws.Range('D13')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.566.22'
$c.Style = $origStyle
$ws.Range('E13').Value = '  -4.28%  '
$c = $ws.Range('D14')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '3.86'
$c.Style = $origStyle
$ws.Range('E14').Value = '  -3.91%  '
$ws.Range('E15').Value = '  -6.06%  '
$c = $ws.Range('D16')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '27.642.98'
$c.Style = $origStyle
$ws.Range('E16').Value = '  -1.10%  '
$c = $ws.Range('D17')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '63.25'
$c.Style = $origStyle
$ws.Range('E17').Value = '  -3.14%  '
$c = $ws.Range('D18')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '219.76'
$c.Style = $origStyle
$ws.Range('E18').Value = '  -4.32%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$c = $ws.Range('D19')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.0₃0692'
$c.Style = $origStyle
$ws.Range('E19').Value = '  -3.87%  '
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c = $ws.Range('D20')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '7.33'
$c.Style = $origStyle
$ws.Range('E20').Value = '  -6.29%  '
$ws.Range('E21').Value = '  +0.37%  '
$c = $ws.Range('D22')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '4.14'
$c.Style = $origStyle
$ws.Range('E22').Value = '  -5.21%  '
$ws.Range('E23').Value = '  -6.27%  '
$ws.Range('E24').Value = '  -5.02%  '
$c = $ws.Range('D25')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '154.08'
$c.Style = $origStyle
$ws.Range('E25').Value = '  -1.33%  '
$c = $ws.Range('D26')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '6.77'
$c.Style = $origStyle
$ws.Range('E26').Value = '  -2.91%  '
$ws.Range('E27').Value = '  +0.39%  '
$c = $ws.Range('D28')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '15.13'
$c.Style = $origStyle
$ws.Range('E28').Value = '  -2.85%  '
$ws.Range('E29').Value = '  -4.18%  '
$c = $ws.Range('D30')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.16'
$c.Style = $origStyle
$ws.Range('E30').Value = '  -2.16%  '
$ws.Range('E31').Value = '  -3.54%  '
$c = $ws.Range('D32')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '3.23'
$c.Style = $origStyle
$ws.Range('E32').Value = '  -5.42%  '
$c = $ws.Range('D33')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.386.32'
$c.Style = $origStyle
$ws.Range('E33').Value = '  -1.11%  '
$ws.Range('E34').Value = '  -5.68%  '
$ws.Range('E35').Value = '  -5.57%  '
$ws.Range('E36').Value = '  -4.54%  '
$ws.Range('E37').Value = '  -0.79%  '
$ws.Range('E38').Value = '  -3.62%  '
$ws.Range('E39').Value = '  -3.60%  '
$ws.Range('E40').Value = '  -4.08%  '
$ws.Range('E42').Value = '  -3.60%  '
$ws.Range('E43').Value = '  +1.40%  '
$c = $ws.Range('D44')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '63.50'
$c.Style = $origStyle
$ws.Range('E44').Value = '  -3.91%  '
$ws.Range('E45').Value = '  -4.69%  '
$c = $ws.Range('D46')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '5.22'
$c.Style = $origStyle
$ws.Range('E46').Value = '  -4.42%  '
$c = $ws.Range('D47')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.719.40'
$c.Style = $origStyle
$ws.Range('E47').Value = '  -3.15%  '
$c = $ws.Range('D48')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '88.02'
$c.Style = $origStyle
$ws.Range('E48').Value = '  -0.64%  '
$ws.Range('E49').Value = '  -2.18%  '
$c = $ws.Range('D50')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.0973'
$c.Style = $origStyle
$ws.Range('E50').Value = '  -4.97%  '
$ws.Range('E51').Value = '  -1.34%  '
